$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# fix(FN-3460): fix invalid facility utilisation values in e2e report fixtures

# Facility limit (column E) corrections - rows 2-6
$ws.Range("E2").Value = 800000
$ws.Range("E3").Value = 800000
$ws.Range("E4").Value = 800000
$ws.Range("E5").Value = 800000
$ws.Range("E6").Value = 800000

# Row 5: fees paid / total fees accrued corrections
$ws.Range("G5").Value = 456
$ws.Range("H5").Value = 3938753.8

# Row 6: fees paid correction
$ws.Range("G6").Value = 761579.37

# Column G now shares the same (bestFit) width as columns E:F since its
# widest value is now "761579.37", matching the rest of the block.
$ws.Range("G1").EntireColumn.ColumnWidth = 15.5

# Reflect the reviewed/selected range left active after the edits.
[void]$ws.Range("E2:H6").Select()
